$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 602 - this pushes the existing row 602 (and
# everything after it) down by one, turning the old row 602 into row 603,
# old row 603 into row 604, etc.
$ws.Rows("602:602").Insert()

# The newly inserted row 602 is a blank row. Populate it as a copy of what
# is now row 603 (the data that used to live in row 602), then overwrite
# the handful of cells (date, min/max/avg price, $/kg) that differ for
# this new weekly record.
$ws.Range("A602:T602").Value2 = $ws.Range("A603:T603").Value2

$ws.Range("D602").Value = 45166
$ws.Range("N602").Value = 14000
$ws.Range("O602").Value = 15000
$ws.Range("P602").Value = 14500
$ws.Range("S602").Value = 1036
